# Update from MV -datos-: add a new "Agosto.2021" quarterly column at the
# right edge of the table, carrying forward the last known value for each
# series (same pattern already used for every prior "no data yet" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$firstCol  = $usedRange.Column
$lastRow   = $firstRow + $usedRange.Rows.Count - 1
$srcCol    = $firstCol + $usedRange.Columns.Count - 1
$newCol    = $srcCol + 1

# Header for the new column: copy the formatting of the previous header
# cell, then set the new period label.
$headerSrc = $ws.Cells.Item($firstRow, $srcCol)
$headerDst = $ws.Cells.Item($firstRow, $newCol)
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122) # xlPasteFormats
$headerDst.Value = "Agosto.2021"

# Every data row simply repeats the last available value into the new
# column (the same way the last existing column already repeats the one
# before it).
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newCol).Value = $ws.Cells.Item($r, $srcCol).Value2
}
